$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.401.59'
$ws.Range("E2").Value = '  +1.55%  '
$ws.Range("D3").Value = '2.276.37'
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.66'
$ws.Range("E5").Value = '  +1.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.55'
$ws.Range("E6").Value = '  +6.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.530'
$ws.Range("E7").Value = '  -0.38%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("E9").Value = '  +2.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.83'
$ws.Range("E10").Value = '  +10.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0797'
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("E12").Value = '  -1.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.70'
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("D14").Value = '2.628.27'
$ws.Range("E14").Value = '  +0.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.43'
$ws.Range("E15").Value = '  +1.36%  '
$ws.Range("D16").Value = '2.260.58'
$ws.Range("E16").Value = '  -0.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.794'
$ws.Range("E17").Value = '  +3.09%  '
$ws.Range("D18").Value = '42.265.51'
$ws.Range("E18").Value = '  +1.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.55'
$ws.Range("E19").Value = '  +1.17%  '
$ws.Range("D20").Value = '0.0₃0911'
$ws.Range("E20").Value = '  +0.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.98'
$ws.Range("E21").Value = '  +0.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.61'
$ws.Range("E22").Value = '  +0.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.10'
$ws.Range("E23").Value = '  +0.59%  '
$ws.Range("E24").Value = '  +0.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.94'
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.82'
$ws.Range("E27").Value = '  -0.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.57'
$ws.Range("E28").Value = '  +7.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.50'
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("E30").Value = '  +1.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '159.03'
$ws.Range("E31").Value = '  -0.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.23'
$ws.Range("E32").Value = '  -0.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.14'
$ws.Range("E34").Value = '  +4.25%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0744'
$ws.Range("E35").Value = '  +0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.106'
$ws.Range("E38").Value = '  +1.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.84'
$ws.Range("E39").Value = '  +2.82%  '
$ws.Range("E40").Value = '  -1.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.10'
$ws.Range("E41").Value = '  +4.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.40'
$ws.Range("E42").Value = '  +13.32%  '
$ws.Range("D43").Value = '1.998.49'
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.10'
$ws.Range("E44").Value = '  -0.85%  '
$ws.Range("E45").Value = '  +2.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.00'
$ws.Range("E46").Value = '  +4.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.97'
$ws.Range("E47").Value = '  -3.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.02'
$ws.Range("E48").Value = '  +1.92%  '
$ws.Range("E49").Value = '  +0.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.22'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '91.93'
$ws.Range("E51").Value = '  +1.24%  '

# Row 36/37: Celestia moved above WEMIXToken in ranking
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.95"
$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.39"
$ws.Range("E37").Value = "  +1.24%  "
